$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.936.96'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -0.49%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.351.85'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -0.60%  '

$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.96'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.30%  '

$ws.Range('E6').Value = '  -3.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.90'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -2.79%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.600'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -3.27%  '

$ws.Range('E10').Value = '  -0.68%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.07'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +4.99%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '33.22'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.06%  '

$ws.Range('E13').Value = '  +0.26%  '

$ws.Range('E14').Value = '  -1.86%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.704.87'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -0.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.19'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -3.10%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.904'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -1.67%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.356.28'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -0.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.876.29'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -0.29%  '

$ws.Range('E20').Value = '  +0.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.64'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -0.18%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.14'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +0.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.76'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.53%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.08%  '

$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.80'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +2.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.85'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +1.64%  '

$ws.Range('E27').Value = '  -0.96%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.44'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -3.92%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -2.11%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.12'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +0.58%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.24'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -3.81%  '

$ws.Range('E32').Value = '  -0.54%  '

$ws.Range('E33').Value = '  -2.27%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0748'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -2.00%  '

$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.07'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -5.04%  '

$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.38'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -0.37%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.82'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +1.49%  '

$ws.Range('E38').Value = '  +0.01%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.42'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +0.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0272'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -4.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.41'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +10.90%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.75'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +13.15%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.16'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -1.28%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.22'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -2.39%  '

$ws.Range('E45').Value = '  -3.93%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.201'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -5.57%  '

$ws.Range('E47').Value = '  +0.01%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.45'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -3.87%  '

$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.23'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -2.36%  '

$ws.Range('E50').Value = '  -2.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '98.41'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -2.75%  '
